{"js": "// Update the 25 two-digit \u00f7 one-digit problems laid out in the single\n// table's 5 \"problem\" rows (rows 0, 4, 8, 12, 16 of the 20-row table;\n// the other rows are blank answer-space rows). Each problem cell's run\n// text is replaced in place (via a scoped search + InsertLocation.replace)\n// so the existing run formatting (rFonts TimeNewRoman, sz 30) is kept\n// untouched \u2014 only the literal digits/text change, matching the diff.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// [rowIndex, colIndex, oldText, newText] in document order.\nconst replacements = [\n  [0, 0, \"76\u00f74=\", \"27\u00f73=\"],\n  [0, 1, \"41\u00f74=\", \"60\u00f79=\"],\n  [0, 2, \"57\u00f73=\", \"15\u00f76=\"],\n  [0, 3, \"43\u00f72=\", \"36\u00f72=\"],\n  [0, 4, \"35\u00f74=\", \"78\u00f77=\"],\n\n  [4, 0, \"81\u00f73=\", \"13\u00f73=\"],\n  [4, 1, \"88\u00f78=\", \"59\u00f75=\"],\n  [4, 2, \"99\u00f76=\", \"49\u00f75=\"],\n  [4, 3, \"13\u00f79=\", \"91\u00f78=\"],\n  [4, 4, \"61\u00f76=\", \"75\u00f78=\"],\n\n  [8, 0, \"74\u00f77=\", \"76\u00f78=\"],\n  [8, 1, \"88\u00f74=\", \"14\u00f78=\"],\n  [8, 2, \"19\u00f73=\", \"93\u00f75=\"],\n  [8, 3, \"24\u00f75=\", \"35\u00f78=\"],\n  [8, 4, \"15\u00f73=\", \"32\u00f75=\"],\n\n  [12, 0, \"83\u00f77=\", \"73\u00f74=\"],\n  [12, 1, \"20\u00f74=\", \"59\u00f79=\"],\n  [12, 2, \"30\u00f77=\", \"89\u00f75=\"],\n  [12, 3, \"79\u00f78=\", \"13\u00f76=\"],\n  [12, 4, \"16\u00f76=\", \"90\u00f76=\"],\n\n  [16, 0, \"20\u00f72=\", \"36\u00f75=\"],\n  [16, 1, \"35\u00f78=\", \"40\u00f77=\"],\n  [16, 2, \"22\u00f75=\", \"32\u00f75=\"],\n  [16, 3, \"59\u00f75=\", \"83\u00f79=\"],\n  [16, 4, \"61\u00f74=\", \"34\u00f77=\"],\n];\n\n// Resolve every target range first (scoped to its own cell, so there is\n// no ambiguity even though some old/new values repeat elsewhere in the\n// table), then apply all the text swaps.\nconst ranges = [];\nfor (const [row, col, oldText] of replacements) {\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  ranges.push(results);\n}\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, , , newText] = replacements[i];\n  ranges[i].items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the 25 two-digit \u00f7 one-digit problems laid out in the single\n# table's 5 \"problem\" rows (table rows 1, 5, 9, 13, 17; the other rows\n# are blank answer-space rows). Each cell is addressed directly by\n# (row, col) and its Range.Text is overwritten in place, so there is no\n# cross-cell ambiguity (the Word.Find object in this host ignores the\n# Range it was obtained from and searches/replaces document-wide, which\n# is unsafe here because some old/new values repeat elsewhere in the\n# table) and the existing run formatting (rFonts TimeNewRoman, sz 30)\n# is left untouched (assigning Range.Text reuses the run's rPr).\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$replacements = @(\n  @{ Row = 1; Col = 1; NewText = \"27\u00f73=\" },\n  @{ Row = 1; Col = 2; NewText = \"60\u00f79=\" },\n  @{ Row = 1; Col = 3; NewText = \"15\u00f76=\" },\n  @{ Row = 1; Col = 4; NewText = \"36\u00f72=\" },\n  @{ Row = 1; Col = 5; NewText = \"78\u00f77=\" },\n\n  @{ Row = 5; Col = 1; NewText = \"13\u00f73=\" },\n  @{ Row = 5; Col = 2; NewText = \"59\u00f75=\" },\n  @{ Row = 5; Col = 3; NewText = \"49\u00f75=\" },\n  @{ Row = 5; Col = 4; NewText = \"91\u00f78=\" },\n  @{ Row = 5; Col = 5; NewText = \"75\u00f78=\" },\n\n  @{ Row = 9; Col = 1; NewText = \"76\u00f78=\" },\n  @{ Row = 9; Col = 2; NewText = \"14\u00f78=\" },\n  @{ Row = 9; Col = 3; NewText = \"93\u00f75=\" },\n  @{ Row = 9; Col = 4; NewText = \"35\u00f78=\" },\n  @{ Row = 9; Col = 5; NewText = \"32\u00f75=\" },\n\n  @{ Row = 13; Col = 1; NewText = \"73\u00f74=\" },\n  @{ Row = 13; Col = 2; NewText = \"59\u00f79=\" },\n  @{ Row = 13; Col = 3; NewText = \"89\u00f75=\" },\n  @{ Row = 13; Col = 4; NewText = \"13\u00f76=\" },\n  @{ Row = 13; Col = 5; NewText = \"90\u00f76=\" },\n\n  @{ Row = 17; Col = 1; NewText = \"36\u00f75=\" },\n  @{ Row = 17; Col = 2; NewText = \"40\u00f77=\" },\n  @{ Row = 17; Col = 3; NewText = \"32\u00f75=\" },\n  @{ Row = 17; Col = 4; NewText = \"83\u00f79=\" },\n  @{ Row = 17; Col = 5; NewText = \"34\u00f77=\" }\n)\n\nforeach ($rep in $replacements) {\n  $table.Cell($rep.Row, $rep.Col).Range.Text = $rep.NewText\n}\n"}
